$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header for new "Shadow casting" weighting column (J)
$ws.Range("J1").Value = "Shadow casting"

# Fill the weighting value (0.7) for each of the 14 data rows.
# The existing sheet already stores "0.7" as text in F10/G10/H10 (same style
# as the rest of column J), so copy that cell's value+format down instead of
# assigning .Value directly (which Excel would auto-convert to a number).
for ($r = 2; $r -le 15; $r++) {
    $ws.Range("F10").Copy()
    $ws.Cells.Item($r, 10).PasteSpecial(-4104)
}
$excel.CutCopyMode = $false

# Re-create the CONCATENATE helper column as a shared formula spanning B20:B33,
# matching the way Excel would fill the formula down the range.
$ws.Range("B20:B33").Formula = '=CONCATENATE(B2,",",C2,",",D2,",",E2,",",F2,",",G2,",",H2,",",I2)'

# Update the selection / view state to match the saved workbook
[void]$ws.Range("G20").Select()
